$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 109, pushing existing rows 109-132 down to 110-133.
$ws.Rows.Item(109).Insert()

# Populate the newly inserted row 109 with the new price-record data.
$ws.Cells.Item(109, 1).Value = 3
$ws.Cells.Item(109, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(109, 3).Value = "Coquimbo"
$ws.Cells.Item(109, 4).Value = 44711
$ws.Cells.Item(109, 5).Value = 5
$ws.Cells.Item(109, 6).Value = "Fruta"
$ws.Cells.Item(109, 7).Value = 100107
$ws.Cells.Item(109, 8).Value = "Otros"
$ws.Cells.Item(109, 9).Value = 100107011
$ws.Cells.Item(109, 10).Value = "Tuna"
$ws.Cells.Item(109, 11).Value = "Sin especificar"
$ws.Cells.Item(109, 12).Value = "Primera"
$ws.Cells.Item(109, 13).Value = 68
$ws.Cells.Item(109, 14).Value = 17000
$ws.Cells.Item(109, 15).Value = 17000
$ws.Cells.Item(109, 16).Value = 17000
$ws.Cells.Item(109, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(109, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(109, 19).Value = 850
$ws.Cells.Item(109, 20).Value = 20
